# Refresh market-price derived columns (H:N) for the leves whose Universalis
# quotes moved since the last scheduled run. Only the cells that actually
# changed are touched; cells that became blank are cleared (not zeroed).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: Days of Chunder / Antidote
$ws.Range("H6").Value = 168146.67
$ws.Range("I6").Value = 168146.67
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 504440.01
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -504328.01
$ws.Range("N6").ClearContents()

# Row 39: Riches' Brew / Hi-Potion of Mind
$ws.Range("H39").Value = 360.8
$ws.Range("I39").Value = 202.25
$ws.Range("K39").Value = 606.75
$ws.Range("M39").Value = -310.75

# Row 45: The House Always Wins / Blinding Potion
$ws.Range("H45").Value = 499.92307
$ws.Range("J45").Value = 499.92307
$ws.Range("L45").Value = 1499.76921
$ws.Range("N45").Value = -1883.76921

# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 6662.5
$ws.Range("J76").Value = 6216.6665
$ws.Range("L76").Value = 6216.6665
$ws.Range("N76").Value = -6846.6665

# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 6662.5
$ws.Range("J79").Value = 6216.6665
$ws.Range("L79").Value = 6216.6665
$ws.Range("N79").Value = -8400.666499999999

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 8714065
$ws.Range("J112").Value = 8714065
$ws.Range("L112").Value = 26142195
$ws.Range("N112").Value = -26144411

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 26702.1
$ws.Range("I137").Value = 45201.19
$ws.Range("J137").Value = 6661.4165
$ws.Range("K137").Value = 135603.57
$ws.Range("L137").Value = 19984.2495
$ws.Range("M137").Value = -133053.57
$ws.Range("N137").Value = -25084.2495

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 4436.022
$ws.Range("I61").Value = 1876.1
$ws.Range("K61").Value = 1876.1
$ws.Range("M61").Value = -1664.1

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 106586.69
$ws.Range("I74").Value = 180967.9
$ws.Range("K74").Value = 180967.9
$ws.Range("M74").Value = -180093.9

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 106586.69
$ws.Range("I77").Value = 180967.9
$ws.Range("K77").Value = 904839.5
$ws.Range("M77").Value = -900471.5

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 4436.022
$ws.Range("I136").Value = 1876.1
$ws.Range("K136").Value = 5628.299999999999
$ws.Range("M136").Value = -3078.299999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 138: Bladewinner / Titanium Gold Greatsword
$ws.Range("H138").Value = 69995.336
$ws.Range("J138").Value = 69995.336
$ws.Range("L138").Value = 69995.336
$ws.Range("N138").Value = -80275.336

$ws = $wb.Worksheets.Item("CRP")
# Row 94: Beech, Please / Beech Lumber
$ws.Range("H94").Value = 1150.2307
$ws.Range("I94").Value = 673.375
$ws.Range("J94").Value = 1362.1666
$ws.Range("K94").Value = 673.375
$ws.Range("L94").Value = 1362.1666
$ws.Range("M94").Value = -222.375
$ws.Range("N94").Value = -2264.1666

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 4985.3784
$ws.Range("I132").Value = 3204.074
$ws.Range("J132").Value = 9794.9
$ws.Range("K132").Value = 9612.222
$ws.Range("L132").Value = 29384.7
$ws.Range("M132").Value = -7082.222
$ws.Range("N132").Value = -34444.7

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 4517.912
$ws.Range("J134").Value = 2213.5
$ws.Range("L134").Value = 6640.5
$ws.Range("N134").Value = -11710.5

$ws = $wb.Worksheets.Item("CUL")
# Row 103: West Meats East / Nomad Meat Pie
$ws.Range("H103").Value = 1400.2222
$ws.Range("J103").Value = 1998.4
$ws.Range("L103").Value = 5995.200000000001
$ws.Range("N103").Value = -7753.200000000001

# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 1291.1111
$ws.Range("J113").Value = 1358.8125
$ws.Range("L113").Value = 4076.4375
$ws.Range("N113").Value = -8416.4375

# Row 126: Imperial Palate / Glory Be Soup
$ws.Range("H126").Value = 4515
$ws.Range("I126").Value = 4030
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 12090
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -7150
$ws.Range("N126").Value = -24880

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers / Copper Ingot
$ws.Range("H2").Value = 94.666664
$ws.Range("I2").Value = 113.14286
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 113.14286
$ws.Range("L2").Value = 30
$ws.Range("M2").Value = -0.1428599999999989
$ws.Range("N2").Value = -256

# Row 7: Water of Life / Copper Rings
$ws.Range("H7").Value = 10002.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 10002.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 10002.5
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -10226.5

# Row 8: Gods of Small Things / Copper Earrings
$ws.Range("H8").Value = 10002.5
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 10002.5
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 10002.5
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -10280.5

# Row 15: The Tusk at Hand / Fang Earrings
$ws.Range("H15").Value = 29990
$ws.Range("J15").Value = 29990
$ws.Range("L15").Value = 29990
$ws.Range("N15").Value = -30566

# Row 49: Faith and Fashion / Mythril Earrings
$ws.Range("H49").Value = 30005
$ws.Range("J49").Value = 30005
$ws.Range("L49").Value = 30005
$ws.Range("N49").Value = -30373

# Row 81: The Grander Temple / Dragon Fang Earrings
$ws.Range("H81").Value = 29990
$ws.Range("J81").Value = 29990
$ws.Range("L81").Value = 29990
$ws.Range("N81").Value = -31986

# Row 84: Man with a Dragon Earring (L) / Dragon Fang Earrings
$ws.Range("H84").Value = 29990
$ws.Range("J84").Value = 29990
$ws.Range("L84").Value = 89970
$ws.Range("N84").Value = -99954

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 8187.174
$ws.Range("I132").Value = 8024
$ws.Range("J132").Value = 8649.5
$ws.Range("K132").Value = 24072
$ws.Range("L132").Value = 25948.5
$ws.Range("M132").Value = -21542
$ws.Range("N132").Value = -31008.5

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 2520.762
$ws.Range("I136").Value = 1061.7188
$ws.Range("J136").Value = 7189.7
$ws.Range("K136").Value = 3185.1564
$ws.Range("L136").Value = 21569.1
$ws.Range("M136").Value = -635.1564000000003
$ws.Range("N136").Value = -26669.1

# Row 139: Giving Gatherers Their Gear / Gomphotherium Doublet of Gathering
$ws.Range("H139").Value = 48329.332
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 48329.332
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 48329.332
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -58609.332

$ws = $wb.Worksheets.Item("WVR")
# Row 29: Getting Handsy / Cotton Dress Gloves
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# Row 41: Half Is the New Double / Linen Halfgloves
$ws.Range("H41").Value = 8742.5
$ws.Range("J41").Value = 8991.714
$ws.Range("L41").Value = 8991.714
$ws.Range("N41").Value = -9771.714

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 271226.8
$ws.Range("I136").Value = 303992.53
$ws.Range("J136").Value = 909.5
$ws.Range("K136").Value = 911977.5900000001
$ws.Range("L136").Value = 2728.5
$ws.Range("M136").Value = -909427.5900000001
$ws.Range("N136").Value = -7828.5
